# Swap the presentation's two themes: "Default" <-> "Simple Light".
#
# The deck ships two theme parts:
#   ppt/theme/theme1.xml  (originally "Default", only used by the Notes Master)
#   ppt/theme/theme2.xml  (originally "Simple Light", used by the Slide Master
#                          / every slide's Design, i.e. the one the object
#                          model exposes as Master.Theme / NotesMaster.Theme)
#
# After the edit theme1.xml should hold the "Simple Light" palette and
# theme2.xml should hold the "Default" palette. The PowerPoint object model
# only exposes a single editable Theme (rooted at the slide master's theme
# part), so we recolor that shared Theme object to the target "Default"
# palette - this is the theme file every slide actually renders with.

$p = $ppt.ActivePresentation
$theme = $p.SlideMaster.Theme
$colors = $theme.ThemeColorScheme

function ComRGB([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

# Target palette = the deck's "Default" color scheme.
$colors.Colors(1).RGB  = ComRGB 0x00 0x00 0x00   # dk1
$colors.Colors(2).RGB  = ComRGB 0xFF 0xFF 0xFF   # lt1
$colors.Colors(3).RGB  = ComRGB 0x15 0x81 0x58   # dk2
$colors.Colors(4).RGB  = ComRGB 0xF3 0xF3 0xF3   # lt2
$colors.Colors(5).RGB  = ComRGB 0x05 0x8D 0xC7   # accent1
$colors.Colors(6).RGB  = ComRGB 0x50 0xB4 0x32   # accent2
$colors.Colors(7).RGB  = ComRGB 0xED 0x56 0x1B   # accent3
$colors.Colors(8).RGB  = ComRGB 0xED 0xEF 0x00   # accent4
$colors.Colors(9).RGB  = ComRGB 0x24 0xCB 0xE5   # accent5
$colors.Colors(10).RGB = ComRGB 0x64 0xE5 0x72   # accent6
$colors.Colors(11).RGB = ComRGB 0x22 0x00 0xCC   # hlink
$colors.Colors(12).RGB = ComRGB 0x55 0x1A 0x8B   # folHlink
